$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1670.5714
$ws.Range("I19").Value = 1574.25
$ws.Range("J19").Value = 1799
$ws.Range("K19").Value = 1574.25
$ws.Range("L19").Value = 1799
$ws.Range("M19").Value = -1399.25
$ws.Range("N19").Value = -2149

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1988.3
$ws.Range("I28").Value = 1752.7693
$ws.Range("J28").Value = 2425.7144
$ws.Range("K28").Value = 1752.7693
$ws.Range("L28").Value = 2425.7144
$ws.Range("M28").Value = -1267.7693
$ws.Range("N28").Value = -3395.7144

# ALC row 63
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 79103.664
$ws.Range("I63").Value = 65000
$ws.Range("J63").Value = 86155.5
$ws.Range("K63").Value = 65000
$ws.Range("L63").Value = 86155.5
$ws.Range("M63").Value = -64376
$ws.Range("N63").Value = -87403.5

# ALC row 66
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H66").Value = 79103.664
$ws.Range("I66").Value = 65000
$ws.Range("J66").Value = 86155.5
$ws.Range("K66").Value = 195000
$ws.Range("L66").Value = 258466.5
$ws.Range("M66").Value = -191880
$ws.Range("N66").Value = -264706.5

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2265.1428
$ws.Range("J88").Value = 2226.0833
$ws.Range("L88").Value = 2226.0833
$ws.Range("N88").Value = -3038.0833

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 2265.1428
$ws.Range("J91").Value = 2226.0833
$ws.Range("L91").Value = 2226.0833
$ws.Range("N91").Value = -5034.0833

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1069
$ws.Range("I92").Value = 1069
$ws.Range("K92").Value = 1069
$ws.Range("M92").Value = 179

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1329.081
$ws.Range("J129").Value = 1749.9584
$ws.Range("L129").Value = 5249.8752
$ws.Range("N129").Value = -15249.8752

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 298808.12
$ws.Range("I138").Value = 2468
$ws.Range("J138").Value = 338320.12
$ws.Range("K138").Value = 7404
$ws.Range("L138").Value = 1014960.36
$ws.Range("M138").Value = -2264
$ws.Range("N138").Value = -1025240.36

# ARM row 75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 30000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 30000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 30000
$ws.Range("M75").Value = $null
$ws.Range("N75").Value = -31748

# ARM row 78
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H78").Value = 30000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 30000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 90000
$ws.Range("M78").Value = $null
$ws.Range("N78").Value = -98736

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 5962.4165
$ws.Range("I97").Value = 1745.8462
$ws.Range("J97").Value = 10945.637
$ws.Range("K97").Value = 1745.8462
$ws.Range("L97").Value = 10945.637
$ws.Range("M97").Value = -1249.8462
$ws.Range("N97").Value = -11937.637

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1981.8064
$ws.Range("I122").Value = 1905.8
$ws.Range("K122").Value = 5717.4
$ws.Range("M122").Value = -3267.4

# ARM row 131
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2206.1904
$ws.Range("J107").Value = 599.75
$ws.Range("L107").Value = 599.75
$ws.Range("N107").Value = -4439.75

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3025.587
$ws.Range("I134").Value = 3048.3777
$ws.Range("K134").Value = 9145.133099999999
$ws.Range("M134").Value = -6610.133099999999

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2384.818
$ws.Range("J31").Value = 6439.6
$ws.Range("L31").Value = 6439.6
$ws.Range("N31").Value = -7029.6

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2384.818
$ws.Range("J34").Value = 6439.6
$ws.Range("L34").Value = 6439.6
$ws.Range("N34").Value = -6843.6

# CRP row 53
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 59999
$ws.Range("J53").Value = 59999
$ws.Range("L53").Value = 59999
$ws.Range("N53").Value = -61213

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1487
$ws.Range("I58").Value = 1217.4445
$ws.Range("J58").Value = 2700
$ws.Range("K58").Value = 1217.4445
$ws.Range("L58").Value = 2700
$ws.Range("M58").Value = -1014.4445
$ws.Range("N58").Value = -3106

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 437.55554
$ws.Range("I107").Value = 419.75
$ws.Range("J107").Value = 451.8
$ws.Range("K107").Value = 419.75
$ws.Range("L107").Value = 451.8
$ws.Range("M107").Value = 1500.25
$ws.Range("N107").Value = -4291.8

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1111.5
$ws.Range("I134").Value = 957.2222
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 2871.6666
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -336.6666
$ws.Range("N134").Value = -12570

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1487
$ws.Range("I136").Value = 1217.4445
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 3652.3335
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -1102.3335
$ws.Range("N136").Value = -13200

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1150
$ws.Range("I99").Value = 1275
$ws.Range("J99").Value = 900
$ws.Range("K99").Value = 3825
$ws.Range("L99").Value = 2700
$ws.Range("M99").Value = -1579
$ws.Range("N99").Value = -7192

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2837.4285
$ws.Range("I132").Value = 1666
$ws.Range("J132").Value = 3306
$ws.Range("K132").Value = 14994
$ws.Range("L132").Value = 29754
$ws.Range("M132").Value = -12464
$ws.Range("N132").Value = -34814

# GSM row 13
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1723.8334
$ws.Range("J13").Value = 4444
$ws.Range("L13").Value = 4444
$ws.Range("N13").Value = -4722

# GSM row 17
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 800
$ws.Range("J17").Value = 800
$ws.Range("L17").Value = 800
$ws.Range("N17").Value = -1136

# GSM row 34
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 46494.8
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 46494.8
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 46494.8
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -47030.8

# GSM row 76
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 46494.8
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 46494.8
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 46494.8
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -47124.8

# GSM row 79
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H79").Value = 46494.8
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 46494.8
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 46494.8
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -48678.8

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2749.25
$ws.Range("I113").Value = 2499
$ws.Range("J113").Value = 2999.5
$ws.Range("K113").Value = 2499
$ws.Range("L113").Value = 2999.5
$ws.Range("M113").Value = -329
$ws.Range("N113").Value = -7339.5

# GSM row 124
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 22475.455
$ws.Range("I126").Value = 4216.5713
$ws.Range("K126").Value = 12649.7139
$ws.Range("M126").Value = -10179.7139

# LTW row 10
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 2212
$ws.Range("I10").Value = 2212
$ws.Range("K10").Value = 2212
$ws.Range("M10").Value = -2072

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 18705.5
$ws.Range("I61").Value = 15207.692
$ws.Range("J61").Value = 27799.8
$ws.Range("K61").Value = 15207.692
$ws.Range("L61").Value = 27799.8
$ws.Range("M61").Value = -15005.692
$ws.Range("N61").Value = -28203.8

# LTW row 63
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 35881.707
$ws.Range("J63").Value = 35881.707
$ws.Range("L63").Value = 35881.707
$ws.Range("N63").Value = -37379.707

# LTW row 66
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 35881.707
$ws.Range("J66").Value = 35881.707
$ws.Range("L66").Value = 107645.121
$ws.Range("N66").Value = -115133.121

# LTW row 75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = $null

# LTW row 76
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 21415.25
$ws.Range("J76").Value = 23800
$ws.Range("L76").Value = 23800
$ws.Range("N76").Value = -24476

# LTW row 78
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = $null

# LTW row 79
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H79").Value = 21415.25
$ws.Range("J79").Value = 23800
$ws.Range("L79").Value = 23800
$ws.Range("N79").Value = -26140

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 18705.5
$ws.Range("I113").Value = 15207.692
$ws.Range("J113").Value = 27799.8
$ws.Range("K113").Value = 15207.692
$ws.Range("L113").Value = 27799.8
$ws.Range("M113").Value = -13037.692
$ws.Range("N113").Value = -32139.8

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4767.561
$ws.Range("I136").Value = 4735.3887
$ws.Range("J136").Value = 4999.2
$ws.Range("K136").Value = 14206.1661
$ws.Range("L136").Value = 14997.6
$ws.Range("M136").Value = -11656.1661
$ws.Range("N136").Value = -20097.6

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3648.8333
$ws.Range("I62").Value = 4165.3335
$ws.Range("K62").Value = 4165.3335
$ws.Range("M62").Value = -3541.3335

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3648.8333
$ws.Range("I65").Value = 4165.3335
$ws.Range("K65").Value = 20826.6675
$ws.Range("M65").Value = -17706.6675

# WVR row 76
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 12000
$ws.Range("I76").Value = 11000
$ws.Range("J76").Value = 13000
$ws.Range("K76").Value = 11000
$ws.Range("L76").Value = 13000
$ws.Range("M76").Value = -10685
$ws.Range("N76").Value = -13630

# WVR row 79
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H79").Value = 12000
$ws.Range("I79").Value = 11000
$ws.Range("J79").Value = 13000
$ws.Range("K79").Value = 11000
$ws.Range("L79").Value = 13000
$ws.Range("M79").Value = -9908
$ws.Range("N79").Value = -15184

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5813.707
$ws.Range("I136").Value = 5076.7075
$ws.Range("J136").Value = 7591.1763
$ws.Range("K136").Value = 15230.1225
$ws.Range("L136").Value = 22773.5289
$ws.Range("M136").Value = -12680.1225
$ws.Range("N136").Value = -27873.5289
